$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.627.55"
$ws.Range("E2").Value = "  +0.61%  "
$ws.Range("D3").Value = "2.540.39"
$ws.Range("E3").Value = "  +2.15%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'524.82"
$ws.Range("E5").Value = "  +0.67%  "
$ws.Range("D6").Value = "'133.51"
$ws.Range("E6").Value = "  -0.98%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  +0.24%  "
$ws.Range("D8").Value = "'0.566"
$ws.Range("E8").Value = "  +1.16%  "
$ws.Range("D9").Value = "2.538.34"
$ws.Range("E9").Value = "  +1.37%  "
$ws.Range("D10").Value = "'0.0983"
$ws.Range("E10").Value = "  -0.73%  "
$ws.Range("E11").Value = "  -1.34%  "
$ws.Range("D12").Value = "'5.18"
$ws.Range("E12").Value = "  -2.42%  "
$ws.Range("E13").Value = "  -2.21%  "
$ws.Range("D14").Value = "2.984.90"
$ws.Range("E14").Value = "  +1.83%  "
$ws.Range("D15").Value = "58.565.40"
$ws.Range("E15").Value = "  +0.60%  "
$ws.Range("D16").Value = "'22.27"
$ws.Range("E16").Value = "  +0.54%  "
$ws.Range("E17").Value = "  -0.25%  "
$ws.Range("D18").Value = "2.535.26"
$ws.Range("E18").Value = "  +1.65%  "
$ws.Range("D19").Value = "'10.70"
$ws.Range("E19").Value = "  +0.04%  "
$ws.Range("D20").Value = "'323.41"
$ws.Range("E20").Value = "  +0.42%  "
$ws.Range("D21").Value = "'4.18"
$ws.Range("E21").Value = "  -0.22%  "
$ws.Range("D22").Value = "'6.16"
$ws.Range("E22").Value = "  +6.72%  "
$ws.Range("E23").Value = "  +0.32%  "
$ws.Range("D24").Value = "'65.23"
$ws.Range("E24").Value = "  +0.79%  "
$ws.Range("E25").Value = "  -1.20%  "
$ws.Range("D26").Value = "'0.999"
$ws.Range("E26").Value = "  +0.38%  "
$ws.Range("E27").Value = "  -0.88%  "
$ws.Range("E28").Value = "  -0.40%  "
$ws.Range("D29").Value = "0.0₃0753"
$ws.Range("E29").Value = "  +0.08%  "
$ws.Range("E30").Value = "  +1.41%  "
$ws.Range("D31").Value = "'168.36"
$ws.Range("E31").Value = "  -0.72%  "
$ws.Range("E32").Value = "  +0.59%  "
$ws.Range("D34").Value = "'0.999"
$ws.Range("E34").Value = "  +0.04%  "
$ws.Range("E35").Value = "  +0.36%  "
$ws.Range("D36").Value = "'18.28"
$ws.Range("E36").Value = "  +0.89%  "
$ws.Range("E37").Value = "  -4.82%  "
$ws.Range("E38").Value = "  -2.28%  "
$ws.Range("E39").Value = "  +0.95%  "
$ws.Range("D40").Value = "'36.50"
$ws.Range("E40").Value = "  -0.46%  "
$ws.Range("D41").Value = "'0.776"
$ws.Range("E41").Value = "  -3.10%  "
$ws.Range("D42").Value = "'278.89"
$ws.Range("E42").Value = "  +0.79%  "
$ws.Range("E43").Value = "  +0.19%  "
$ws.Range("E44").Value = "  -0.84%  "
$ws.Range("D45").Value = "'0.605"
$ws.Range("E45").Value = "  +0.82%  "
$ws.Range("D46").Value = "'130.00"
$ws.Range("E46").Value = "  +4.59%  "
$ws.Range("E47").Value = "  +0.76%  "
$ws.Range("E48").Value = "  +1.79%  "
$ws.Range("D49").Value = "'17.81"
$ws.Range("E49").Value = "  +0.18%  "
$ws.Range("E50").Value = "  +0.06%  "
$ws.Range("D51").Value = "'17.08"
$ws.Range("E51").Value = "  -0.64%  "
